$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 234, shifting the existing rows 234-250 down to 235-251.
$ws.Rows.Item(234).Insert()

# Populate the newly inserted row 234 with the new "Especial" quality record.
$ws.Cells.Item(234, 1).Value = 10
$ws.Cells.Item(234, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(234, 3).Value = "La Araucanía"
$ws.Cells.Item(234, 4).Value = 44714
$ws.Cells.Item(234, 5).Value = 9
$ws.Cells.Item(234, 6).Value = "Fruta"
$ws.Cells.Item(234, 7).Value = 100101
$ws.Cells.Item(234, 8).Value = "Berries"
$ws.Cells.Item(234, 9).Value = 100112025
$ws.Cells.Item(234, 10).Value = "Frutilla"
$ws.Cells.Item(234, 11).Value = "Sin especificar"
$ws.Cells.Item(234, 12).Value = "Especial"
$ws.Cells.Item(234, 13).Value = 600
$ws.Cells.Item(234, 14).Value = 18000
$ws.Cells.Item(234, 15).Value = 18000
$ws.Cells.Item(234, 16).Value = 18000
$ws.Cells.Item(234, 17).Value = "$/bandeja 7 kilos"
$ws.Cells.Item(234, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(234, 19).Value = 2571
$ws.Cells.Item(234, 20).Value = 7
